$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Set Runmode column (C) to "Y" for every suite row (rows 2-7),
# so all suites -- including Suite E -- are flagged to run.
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Move the active selection to D8 (as left by the author after editing).
$ws.Range("D8").Select()
